$d = $word.ActiveDocument

# Step 1: delete the "GitHub; " run entirely
$d.Content.Find.Execute("GitHub; ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)

# Step 2: change "GitLab" to "Git"
$d.Content.Find.Execute("GitLab", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Git", 2)
